# Weekly update: a new price-report row is inserted for
# "Terminal Hortofrutícola Agro Chillán" / Pera / Packham's Triumph.
# The existing row 71 (Especial / 60 / 11000 / 11000 / 11000 / 688) is the
# most recent sample; a new sample for a later date (44665) replaces it at
# row 71, and the previous sample's data is pushed down into a newly
# inserted row 72 (keeping its original date 44351). Every row below that
# shifts down by one, which also grows the sheet's used range from
# A1:T181 to A1:T182.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new entire row above the current row 72; this shifts rows
# 72..181 down to 73..182 and extends the dimension accordingly.
$ws.Rows.Item(72).Insert()

# The newly inserted row 72 becomes a copy of what row 71 used to contain
# (same variety/quality/volume/price data), keeping the older date.
$ws.Range("A72").Value = 7
$ws.Range("B72").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C72").Value = "Ñuble"
$ws.Range("D72").Value = 44351
$ws.Range("E72").Value = 16
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100104
$ws.Range("H72").Value = "Frutos de pepita"
$ws.Range("I72").Value = 100104005
$ws.Range("J72").Value = "Pera"
$ws.Range("K72").Value = "Packham's Triumph"
$ws.Range("L72").Value = "Especial"
$ws.Range("M72").Value = 60
$ws.Range("N72").Value = 11000
$ws.Range("O72").Value = 11000
$ws.Range("P72").Value = 11000
$ws.Range("Q72").Value = "$/caja 16 kilos empedrada"
$ws.Range("R72").Value = "Provincia de Curicó"
$ws.Range("S72").Value = 688
$ws.Range("T72").Value = 16

# Row 71 now records the newer sample: only the date actually changes
# (quality/volume/prices stay the same as before).
$ws.Range("D71").Value = 44665
